# Update "想去人数" (column F) figures to the freshly scraped counts
# committed as "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 701   # F2: 694 -> 701
$ws.Cells.Item(3, 6).Value = 13502   # F3: 13444 -> 13502
$ws.Cells.Item(4, 6).Value = 13686   # F4: 13612 -> 13686
$ws.Cells.Item(6, 6).Value = 1349   # F6: 1350 -> 1349
$ws.Cells.Item(7, 6).Value = 1356   # F7: 1350 -> 1356
$ws.Cells.Item(8, 6).Value = 5714   # F8: 5688 -> 5714
$ws.Cells.Item(9, 6).Value = 955   # F9: 954 -> 955
$ws.Cells.Item(12, 6).Value = 378   # F12: 377 -> 378
$ws.Cells.Item(13, 6).Value = 202   # F13: 201 -> 202
$ws.Cells.Item(14, 6).Value = 1494   # F14: 1488 -> 1494
$ws.Cells.Item(15, 6).Value = 404   # F15: 400 -> 404
$ws.Cells.Item(16, 6).Value = 2109   # F16: 2097 -> 2109
$ws.Cells.Item(17, 6).Value = 1136   # F17: 1126 -> 1136
$ws.Cells.Item(18, 6).Value = 1731   # F18: 1714 -> 1731
$ws.Cells.Item(19, 6).Value = 901   # F19: 900 -> 901
$ws.Cells.Item(21, 6).Value = 2225   # F21: 2219 -> 2225
$ws.Cells.Item(22, 6).Value = 539   # F22: 537 -> 539
$ws.Cells.Item(23, 6).Value = 772   # F23: 764 -> 772
$ws.Cells.Item(24, 6).Value = 3174   # F24: 3160 -> 3174
$ws.Cells.Item(25, 6).Value = 310   # F25: 309 -> 310
$ws.Cells.Item(26, 6).Value = 285   # F26: 286 -> 285
$ws.Cells.Item(27, 6).Value = 2249   # F27: 2229 -> 2249
$ws.Cells.Item(28, 6).Value = 47   # F28: 42 -> 47
$ws.Cells.Item(31, 6).Value = 1736   # F31: 1733 -> 1736
$ws.Cells.Item(32, 6).Value = 1047   # F32: 1044 -> 1047
$ws.Cells.Item(33, 6).Value = 1265   # F33: 1242 -> 1265
$ws.Cells.Item(34, 6).Value = 84   # F34: 83 -> 84
$ws.Cells.Item(35, 6).Value = 123   # F35: 121 -> 123
$ws.Cells.Item(36, 6).Value = 4512   # F36: 4489 -> 4512
$ws.Cells.Item(37, 6).Value = 4625   # F37: 4604 -> 4625
$ws.Cells.Item(39, 6).Value = 145   # F39: 144 -> 145
$ws.Cells.Item(40, 6).Value = 641   # F40: 640 -> 641
$ws.Cells.Item(41, 6).Value = 665   # F41: 663 -> 665
$ws.Cells.Item(42, 6).Value = 3223   # F42: 3212 -> 3223
$ws.Cells.Item(45, 6).Value = 319   # F45: 317 -> 319
$ws.Cells.Item(46, 6).Value = 71   # F46: 66 -> 71
$ws.Cells.Item(47, 6).Value = 52   # F47: 49 -> 52
$ws.Cells.Item(48, 6).Value = 4370   # F48: 4358 -> 4370
$ws.Cells.Item(49, 6).Value = 247   # F49: 241 -> 247

# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 104   # F4: 102 -> 104
$ws.Cells.Item(6, 6).Value = 1   # F6: 0 -> 1
$ws.Cells.Item(20, 6).Value = 10   # F20: 9 -> 10

# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 7276   # F2: 7242 -> 7276
$ws.Cells.Item(3, 6).Value = 189   # F3: 179 -> 189
$ws.Cells.Item(4, 6).Value = 580   # F4: 550 -> 580

# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 701   # F2: 694 -> 701
$ws.Cells.Item(3, 6).Value = 189   # F3: 179 -> 189
$ws.Cells.Item(4, 6).Value = 580   # F4: 550 -> 580
$ws.Cells.Item(6, 6).Value = 13502   # F6: 13445 -> 13502
$ws.Cells.Item(7, 6).Value = 13686   # F7: 13613 -> 13686
$ws.Cells.Item(9, 6).Value = 1349   # F9: 1350 -> 1349
$ws.Cells.Item(10, 6).Value = 5714   # F10: 5689 -> 5714
$ws.Cells.Item(11, 6).Value = 955   # F11: 954 -> 955
$ws.Cells.Item(12, 6).Value = 378   # F12: 377 -> 378
$ws.Cells.Item(13, 6).Value = 104   # F13: 102 -> 104
$ws.Cells.Item(14, 6).Value = 202   # F14: 201 -> 202
$ws.Cells.Item(15, 6).Value = 1   # F15: 0 -> 1
$ws.Cells.Item(16, 6).Value = 1494   # F16: 1488 -> 1494
$ws.Cells.Item(17, 6).Value = 404   # F17: 400 -> 404
$ws.Cells.Item(18, 6).Value = 2109   # F18: 2097 -> 2109
$ws.Cells.Item(19, 6).Value = 1136   # F19: 1126 -> 1136
$ws.Cells.Item(20, 6).Value = 1731   # F20: 1714 -> 1731
$ws.Cells.Item(21, 6).Value = 901   # F21: 900 -> 901
$ws.Cells.Item(22, 6).Value = 539   # F22: 537 -> 539
$ws.Cells.Item(23, 6).Value = 3174   # F23: 3160 -> 3174
$ws.Cells.Item(24, 6).Value = 285   # F24: 286 -> 285
$ws.Cells.Item(25, 6).Value = 47   # F25: 42 -> 47
$ws.Cells.Item(28, 6).Value = 1736   # F28: 1733 -> 1736
$ws.Cells.Item(31, 6).Value = 1265   # F31: 1242 -> 1265
$ws.Cells.Item(33, 6).Value = 84   # F33: 83 -> 84
$ws.Cells.Item(34, 6).Value = 4512   # F34: 4489 -> 4512
$ws.Cells.Item(35, 6).Value = 4625   # F35: 4604 -> 4625
$ws.Cells.Item(38, 6).Value = 145   # F38: 144 -> 145
$ws.Cells.Item(39, 6).Value = 641   # F39: 640 -> 641
$ws.Cells.Item(40, 6).Value = 665   # F40: 663 -> 665
$ws.Cells.Item(41, 6).Value = 3223   # F41: 3212 -> 3223
$ws.Cells.Item(44, 6).Value = 319   # F44: 317 -> 319
$ws.Cells.Item(45, 6).Value = 71   # F45: 66 -> 71
$ws.Cells.Item(46, 6).Value = 52   # F46: 49 -> 52
$ws.Cells.Item(47, 6).Value = 4370   # F47: 4358 -> 4370
$ws.Cells.Item(48, 6).Value = 247   # F48: 241 -> 247
